# fix SRS requirements numbering
#
# The "Req ID" column (B) had a duplicated "R08" in rows 12 and 13, which
# pushed every subsequent requirement ID out of sequence. Renumber rows
# 13-16 so the sequence reads R08, R09, R10, R11, R12 (row 16 previously
# reused "R11" and now becomes the new, correctly numbered "R12").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "R09"
$ws.Range("B14").Value = "R10"
$ws.Range("B15").Value = "R11"
$ws.Range("B16").Value = "R12"

# Restore the view to where the author left it: scrolled down so row 7 is
# at the top, zoomed to 81%, with J13 as the active/selected cell.
$win = $excel.ActiveWindow
$win.Zoom = 81
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("J13").Select()
